$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "See you in [CITY NAME]!" "Vemo-nos em [CITY NAME]!"
Replace-Text "Great news! We will be in [CITY NAME] from" "Temos ótimas notícias! Vamos estar em [CITY NAME] de "
Replace-Text "[DATE] to [DATE] 2023. Our affiliate team, led by [COUNTRY] Country Manager [AFFILIATE MANAGER NAME], look forward to an exclusive one-on-one session with you." "[DATE] a [DATE] 2023. A nossa equipa de afiliados, liderada pelo Gestor de Parcerias de [PAÍS] [NOME DO AFFILIATE MANAGER], terá todo o gosto em ter uma reunião consigo."
Replace-Text "We’d love to hear about your experience with our affiliate programme. If there’s any way we can improve your experience, here’s your chance to tell us." "Gostaríamos de saber mais sobre a sua experiência com o nosso programa de afiliados. Caso haja algo que possamos melhor, esta é a sua oportunidade de falar connosco."
Replace-Text "When?" "Quando?"
Replace-Text "A 1-hour slot between 9:00 AM and 6:00 PM" "Uma sessão de 1 hora, entre as 9:00h e as 18:00h,"
Replace-Text "from [DATE] to [DATE]" "de [DATE] a [DATE]"
Replace-Text "Where?" "Onde?"
Replace-Text "To be confirmed" "A confirmar"
Replace-Text "How to book a slot?" "Como reservar um horário?"
Replace-Text "Pick a date and time, and reply to this email by [DATE]  (first come, first served)" "Escolha uma data e hora e responda a este e-mail até dia [DATA] (por ordem de submissão)"
Replace-Text "You’re welcome to bring along your clients and friends interested in learning more about trading on Deriv." "Pode trazer os seus clientes e amigos interessados em aprender mais sobre como negociar na Deriv."
Replace-Text "We’re grateful for your continuous support and look forward to meeting you!" "Agradecemos o seu apoio contínuo e esperamos vê-lo em breve!"
Replace-Text "P.S. We’re giving out free Deriv merchandise. Don’t miss out!" "P.D. P.S. Vamos entregar brindes da Deriv. Não perca!"
Replace-Text "If you have questions, contact us " "Caso tenha alguma dúvida, entre em contato connosco "
Replace-Text "P.S. Iremos distribuir vários brindes da Deriv gratuitamente. Não perca! " "P.D. Iremos distribuir vários brindes da Deriv gratuitamente. Não perca! "
Replace-Text "P.S. Nous distribuons des produits Deriv gratuits. Ne manquez pas!" "P.D. Nous distribuons des produits Deriv gratuits. Ne manquez pas!"
Replace-Text "P.S. Chúng tôi sẽ tặng các sản phẩm quảng bá miễn phí của Deriv. Đừng bỏ lỡ cơ hội này!" "P.D. Chúng tôi sẽ tặng các sản phẩm quảng bá miễn phí của Deriv. Đừng bỏ lỡ cơ hội này!"
